$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2's value moves to D2, B2 becomes "-", C2 stays "-"
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "[-, 'MCT-3A-Lab. de eletroeletrônica', -]"

# Row 3: D3 gets a new schedule value, F3 becomes "-"
$ws.Range("D3").Value = "[-, 'MCT-3A-Lab. de eletroeletrônica', -]"
$ws.Range("F3").Value = "-"

# Row 6: D6 gets a new schedule value
$ws.Range("D6").Value = "['MCT-2A-Sistemas digitais', -, 'MCT-2A-Sistemas digitais']"

# Row 8: B8 becomes "-"
$ws.Range("B8").Value = "-"

# Row 18: E18 updated value
$ws.Range("E18").Value = "['ELM-2NA-Lab. Circuitos Elétricos', 'ELM-2NA-Lab. Circuitos Elétricos']"

# Row 20: B20 becomes "-"
$ws.Range("B20").Value = "-"
